# Apply the "Merging Excel sheets is now possible" change:
#  - Rewrite the explanatory text block (rows 12-16) with the new,
#    shorter limitations text, styled in a bold accent color.
#  - Add a new "Example #1: merging data with rows" sub-title (row 18).
#  - Add a 2-row x 3-col demo table (headers row 20, TBS tags row 21)
#    with header shading/border and data border styling.
#  - Resize columns B:D and move the selection to B16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Clear out the old "limitations" paragraphs (rows 12-21) so we can
#    rebuild the block from scratch at the new row positions (also
#    drops stray formatting left on untouched gap rows/cells).
# ---------------------------------------------------------------------
$ws.Range("B12:D21").Clear()

# ---------------------------------------------------------------------
# 2. New text block (rows 12-16) - bold, accent-colored paragraphs.
# ---------------------------------------------------------------------
$introText = "You may consider the following before building your own Microsoft Excel template:"
$limitationLines = @(
    $introText,
    "Merging Microsoft Excel templates with OpenTBS has several limitations because of the OpenXML format for Excel.",
    "* Formulas won't work because OpenTBS needs to convert cell positions from aboslute to relative in order to have a constistent merged sheet.",
    "* Formulas may also make troubles because they are saved twice in the sheet:  one for the expression, and one for the instant result.",
    "* Changing picture (using ope=changepic)  because drawing information are saved in another XML sub-file."
)

$r = 12
foreach ($line in $limitationLines) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $line
    $cell.Font.Bold = $true
    $cell.Font.ThemeColor = 6
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. "Example #1" sub-title (row 18) - bold, normal (theme text) color.
# ---------------------------------------------------------------------
$ws.Cells.Item(18, 2).Value2 = "Example #1: merging data with rows"
$ws.Cells.Item(18, 2).Font.Bold = $true
$ws.Cells.Item(18, 2).Font.ThemeColor = 1

# ---------------------------------------------------------------------
# 4. Demo table: header row (20) + TBS-tag data row (21).
# ---------------------------------------------------------------------
$ws.Cells.Item(20, 2).Value2 = "First Name"
$ws.Cells.Item(20, 3).Value2 = "Name"
$ws.Cells.Item(20, 4).Value2 = "Membership number"

$headerRange = $ws.Range("B20:D20")
$headerRange.Interior.ColorIndex = 64
$headerRange.Interior.ThemeColor = 2
$headerRange.Interior.Pattern = 1
$headerRange.Borders.ColorIndex = 64
$headerRange.Borders.LineStyle = 1

$ws.Cells.Item(21, 2).Value2 = "[a.firstname;block=row]"
$ws.Cells.Item(21, 3).Value2 = "[a.name]"
$ws.Cells.Item(21, 4).Value2 = "[a.number]"

$dataRange = $ws.Range("B21:D21")
$dataRange.Borders.ColorIndex = 64
$dataRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 5. Title cell (B2) keeps its text ("OpenTBS demo") but changes style
#    (font index shuffled in styles.xml - same visual bold Times title).
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 2).Value2 = "OpenTBS demo"

# ---------------------------------------------------------------------
# 6. Column widths for the new table (B, C, D).
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 15.8333333333333
$ws.Columns.Item(3).ColumnWidth = 12.1666666666667
$ws.Columns.Item(4).ColumnWidth = 19

# ---------------------------------------------------------------------
# 7. Selection moves to B16 (matches the new last line of the
#    limitations paragraph).
# ---------------------------------------------------------------------
$ws.Range("B16").Select() | Out-Null
